$wb = $excel.ActiveWorkbook

# --- Sheet rename: ENV -> ENV_URLS ---
$wsEnv = $wb.Worksheets.Item(1)
$wsHome = $wb.Worksheets.Item(2)
$wsEnv.Name = "ENV_URLS"

# --- Add two new blank worksheets at the end (Sheet1, Sheet2) ---
$sheet1 = $wb.Worksheets.Add($null, $wsHome)
$sheet2 = $wb.Worksheets.Add($null, $sheet1)

# --- ENV_URLS: add new column A width ---
$wsEnv.Columns.Item(1).ColumnWidth = 23.998

# --- ENV_URLS: new rows of data ---
# Insert strings in an order that reproduces the target shared-string table order.
$wsEnv.Range("A3").Value = "asdfasdfas"
$wsEnv.Range("A14").Value = "asdfadfas"
$wsEnv.Range("A19").Value = "asdfasdfa"
$wsEnv.Range("A6").Value = "asdfasdfas 6"
$wsEnv.Range("A4").Value = "asdfasdfas"
$wsEnv.Range("A5").Value = "asdfasdfas"

# Currency-formatted number
$wsEnv.Range("A7").Value = 123.45
$wsEnv.Range("A7").Style = "Currency"
$wsEnv.Range("A7").NumberFormat = '_-[$$-1009]* #,##0.00_-;\-[$$-1009]* #,##0.00_-;_-[$$-1009]* "-"??_-;_-@_-'

# Plain numbers + SUM formula
$wsEnv.Range("A8").Value = 1
$wsEnv.Range("A9").Value = 2
$wsEnv.Range("A10").Formula = "=SUM(A8:A9)"

# --- Selection / active-sheet bookkeeping ---
$null = $wsHome.Range("B3").Select()
[void]$wsEnv.Activate()
$null = $wsEnv.Range("B7").Select()
